$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" date column (C2:C15) from serial 45170 (2023-09-01)
# to serial 45174 (2023-09-05), keeping the existing date formatting.
$ws.Range("C2:C15").Value = 45174
